# Edit script for LOM3232.xlsx
# Inserts a new row (Docentes responsaveis value row) after row 12,
# pushing subsequent rows down by one, then updates the text content
# of several B/C cells that changed as part of the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 13 (shifts rows 13:23 down to 14:24)
$ws.Rows.Item(13).Insert()

# The insert clones row 12s column-A (bold label) formatting into the new
# A13 cell; row 13 has no label in the target layout, so drop it entirely.
$ws.Range("A13").Clear()

# Give B13 the same "normal / wrap text" formatting used by the rest of column B
# (the freshly inserted row otherwise inherits the bold label style from column A).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Update the cell text that changed in this revision
$ws.Range("B10").Value = 'Transmitir aos alunos o conhecimento básico sobre metrologia óptica ou seja métodos de medição de tamanho e geometria de componentes mecânicos com o emprego de métodos ópticos, com particular ênfase na interferometria a laser.'
$ws.Range("C10").Value = 'Transmitir aos alunos o conhecimento básico sobre metrologia óptica ou seja métodos de medição de tamanho e geometria de componentes mecânicos com o emprego de métodos ópticos, com particular ênfase na interferometria a laser.'
$ws.Range("B13").Value = '5840793 - Sérgio Schneider'
$ws.Range("C13").Value = '5840793 - Sérgio Schneider'
$ws.Range("B14").Value = 'Apresentar as principais técnicas ópticas para a medição de grandezas como comprimento, deslocamento e forma, com ênfase nas técnicas interferométricas a laser.'
$ws.Range("C14").Value = 'Apresentar as principais técnicas ópticas para a medição de grandezas como comprimento, deslocamento e forma, com ênfase nas técnicas interferométricas a laser.'
$ws.Range("B16").Value = 'Teoria eletromagnética da luz: noções de representação matemática da onda de luz e interpretação de fenômenos como polarização, interferência e difração. Refração, reflexão e óptica geométrica: leis de Snell, equações de Fraunhofer, reflexão total e óptica geométrica. Propagação da luz em meios especiais como cristais fibras ópticas. Óptica de Fourier e holografia: transformada de Fourier e a sua aplicação na óptica como caso de filtros especiais e halográfia. Fontes e sensores de luz: definição e descrição de fontes incoerentes e coerentes e descrição de sensores do tipo puntual, de posição e de imagem. Componentes ópticos e ajuste de sistemas ópticos. Medição de comprimento: método como interferometria, franjas de Moirè, métodos para medição de grandes distâncias. Medição de forma: diversos métodos e técnicas para medição de forma geométrica. Medição de deslocamento, deformação e vibração: métodos de medição que empregam a holografia, speckle" e as franjas de Moirè. Medição de velocidade: métodos de medição de velocidade e sensor de fibras ópticas. Inspeção de falhas: métodos para inspeção de falhas geométricas e internas utilizando a difração ou a difusão da luz.'
$ws.Range("C16").Value = 'Teoria eletromagnética da luz: noções de representação matemática da onda de luz e interpretação de fenômenos como polarização, interferência e difração. Refração, reflexão e óptica geométrica: leis de Snell, equações de Fraunhofer, reflexão total e óptica geométrica. Propagação da luz em meios especiais como cristais fibras ópticas. Óptica de Fourier e holografia: transformada de Fourier e a sua aplicação na óptica como caso de filtros especiais e halográfia. Fontes e sensores de luz: definição e descrição de fontes incoerentes e coerentes e descrição de sensores do tipo puntual, de posição e de imagem. Componentes ópticos e ajuste de sistemas ópticos. Medição de comprimento: método como interferometria, franjas de Moirè, métodos para medição de grandes distâncias. Medição de forma: diversos métodos e técnicas para medição de forma geométrica. Medição de deslocamento, deformação e vibração: métodos de medição que empregam a holografia, speckle" e as franjas de Moirè. Medição de velocidade: métodos de medição de velocidade e sensor de fibras ópticas. Inspeção de falhas: métodos para inspeção de falhas geométricas e internas utilizando a difração ou a difusão da luz.'
$ws.Range("B19").Value = 'Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios.'
$ws.Range("C19").Value = 'Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios.'
$ws.Range("B20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'
$ws.Range("C20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("B22").Value = 'YOSHIZAWA, T. Handbook of Optical Metrology, Boca Raton: CRC Press, 2009.
SALEH, B. E. A.; TEICH, M. C. Handbook of Fotonics, Wiley-Interscience, 2007.
JENKINS, F. A.; WHITE, H. E. Fundamentals of Optics, McGraw-Hill, 1981. 
CREATH, H.; WYANT, J. Measurement of ultraprecision components using non-contact interferometry based instrumentation, Ultraprecision in Manufacturing Engineering, Springer Verlag, 1988.'
$ws.Range("C22").Value = 'YOSHIZAWA, T. Handbook of Optical Metrology, Boca Raton: CRC Press, 2009.
SALEH, B. E. A.; TEICH, M. C. Handbook of Fotonics, Wiley-Interscience, 2007.
JENKINS, F. A.; WHITE, H. E. Fundamentals of Optics, McGraw-Hill, 1981. 
CREATH, H.; WYANT, J. Measurement of ultraprecision components using non-contact interferometry based instrumentation, Ultraprecision in Manufacturing Engineering, Springer Verlag, 1988.'

Write-Output "done"
